# "fixed datapath and finish linking in ID"
# Adds the final control-signal row (RegWrite / WB / null / "write to the
# register file") to the ID-stage signal table on Sheet1, and moves the
# window/selection to reflect that new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: signal name, stage it's generated in, value label, description.
$ws.Range("A13").Value = "RegWrite"
$ws.Range("B13").Value = "WB"
$ws.Range("C13").Value = "null"
$ws.Range("D13").Value = "写入寄存器堆"

# Reposition the window / active cell the way Excel leaves things after
# typing the new row (used range grows to A1:E13, selection lands on D13).
$excel.ActiveWindow.Left = 3960
$ws.Range("D13").Select() | Out-Null
